# Edit script implementing:
# - Update "Created" date string (Proposal!U2) from 03/05/20 to 06/17/20
# - Update Contract Tab (Proposal) quarterly-delivery numeric figures for the
#   "multiple spot lengths" recompute (rows 9-12,17-20,24-27,33-36,41-44,48-51)
# - Apply ShrinkToFit formatting to the Flow Chart section headers (B7,B16,B25,B34,B43,B52)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Created" date on the Proposal sheet (drives the Flow Chart's formula too)
# ---------------------------------------------------------------------------
$proposal = $wb.Worksheets.Item("Proposal")
$proposal.Range("U2").Value = "Created 06/17/20"

# ---------------------------------------------------------------------------
# 2. Quarterly delivery numeric recompute on the Proposal sheet
# ---------------------------------------------------------------------------

$ws = $proposal

$ws.Range("J9").Value = 8.275
$ws.Range("N9").Value = 15105.740181268882
$ws.Range("Q9").Value = 31.3
$ws.Range("U9").Value = 3993.6102236421725
$ws.Range("J10").Value = 8.275
$ws.Range("N10").Value = 15105.740181268882
$ws.Range("Q10").Value = 31.3
$ws.Range("U10").Value = 3993.6102236421725
$ws.Range("J11").Value = 16.55
$ws.Range("N11").Value = 15105.740181268882
$ws.Range("Q11").Value = 62.6
$ws.Range("U11").Value = 3993.6102236421725
$ws.Range("J12").Value = 33.1
$ws.Range("N12").Value = 15105.740181268882
$ws.Range("Q12").Value = 125.2
$ws.Range("U12").Value = 3993.6102236421725
$ws.Range("H17").Value = 1.2249285749999999
$ws.Range("J17").Value = 8.175
$ws.Range("L17").Value = 12249.28575
$ws.Range("M17").Value = 10.204676627778074
$ws.Range("N17").Value = 15290.51987767584
$ws.Range("O17").Value = 1.7819580999999998
$ws.Range("Q17").Value = 8.3
$ws.Range("S17").Value = 17819.581
$ws.Range("T17").Value = 7.014755285211252
$ws.Range("U17").Value = 15060.240963855422
$ws.Range("H18").Value = 1.2249285749999999
$ws.Range("J18").Value = 8.175
$ws.Range("L18").Value = 12249.28575
$ws.Range("M18").Value = 10.204676627778074
$ws.Range("N18").Value = 15290.51987767584
$ws.Range("O18").Value = 1.7819580999999998
$ws.Range("Q18").Value = 8.3
$ws.Range("S18").Value = 17819.581
$ws.Range("T18").Value = 7.014755285211252
$ws.Range("U18").Value = 15060.240963855422
$ws.Range("H19").Value = 1.2249285749999999
$ws.Range("J19").Value = 16.35
$ws.Range("L19").Value = 24498.5715
$ws.Range("M19").Value = 10.204676627778074
$ws.Range("N19").Value = 15290.51987767584
$ws.Range("O19").Value = 1.7819580999999998
$ws.Range("Q19").Value = 16.6
$ws.Range("S19").Value = 35639.162
$ws.Range("T19").Value = 7.014755285211252
$ws.Range("U19").Value = 15060.240963855422
$ws.Range("J20").Value = 32.7
$ws.Range("L20").Value = 48997.143
$ws.Range("M20").Value = 10.204676627778074
$ws.Range("N20").Value = 15290.51987767584
$ws.Range("Q20").Value = 33.2
$ws.Range("S20").Value = 71278.324
$ws.Range("T20").Value = 7.014755285211252
$ws.Range("U20").Value = 15060.240963855422
$ws.Range("H24").Value = 0.49617645000000005
$ws.Range("J24").Value = 8.35
$ws.Range("L24").Value = 4961.7645
$ws.Range("M24").Value = 25.19265071931568
$ws.Range("N24").Value = 14970.05988023952
$ws.Range("H25").Value = 0.49617645000000005
$ws.Range("J25").Value = 8.35
$ws.Range("L25").Value = 4961.7645
$ws.Range("M25").Value = 25.19265071931568
$ws.Range("N25").Value = 14970.05988023952
$ws.Range("H26").Value = 0.49617645000000005
$ws.Range("J26").Value = 16.7
$ws.Range("L26").Value = 9923.529
$ws.Range("M26").Value = 25.19265071931568
$ws.Range("N26").Value = 14970.05988023952
$ws.Range("J27").Value = 33.4
$ws.Range("L27").Value = 19847.058
$ws.Range("M27").Value = 25.19265071931568
$ws.Range("N27").Value = 14970.05988023952
$ws.Range("J33").Value = 8.275
$ws.Range("N33").Value = 15105.740181268882
$ws.Range("Q33").Value = 31.3
$ws.Range("U33").Value = 3993.6102236421725
$ws.Range("J34").Value = 8.275
$ws.Range("N34").Value = 15105.740181268882
$ws.Range("Q34").Value = 31.3
$ws.Range("U34").Value = 3993.6102236421725
$ws.Range("J35").Value = 16.55
$ws.Range("N35").Value = 15105.740181268882
$ws.Range("Q35").Value = 62.6
$ws.Range("U35").Value = 3993.6102236421725
$ws.Range("J36").Value = 33.1
$ws.Range("N36").Value = 15105.740181268882
$ws.Range("Q36").Value = 125.2
$ws.Range("U36").Value = 3993.6102236421725
$ws.Range("H41").Value = 1.2249285749999999
$ws.Range("J41").Value = 8.175
$ws.Range("L41").Value = 12249.28575
$ws.Range("M41").Value = 10.204676627778074
$ws.Range("N41").Value = 15290.51987767584
$ws.Range("O41").Value = 1.7819580999999998
$ws.Range("Q41").Value = 8.3
$ws.Range("S41").Value = 17819.581
$ws.Range("T41").Value = 7.014755285211252
$ws.Range("U41").Value = 15060.240963855422
$ws.Range("H42").Value = 1.2249285749999999
$ws.Range("J42").Value = 8.175
$ws.Range("L42").Value = 12249.28575
$ws.Range("M42").Value = 10.204676627778074
$ws.Range("N42").Value = 15290.51987767584
$ws.Range("O42").Value = 1.7819580999999998
$ws.Range("Q42").Value = 8.3
$ws.Range("S42").Value = 17819.581
$ws.Range("T42").Value = 7.014755285211252
$ws.Range("U42").Value = 15060.240963855422
$ws.Range("H43").Value = 1.2249285749999999
$ws.Range("J43").Value = 16.35
$ws.Range("L43").Value = 24498.5715
$ws.Range("M43").Value = 10.204676627778074
$ws.Range("N43").Value = 15290.51987767584
$ws.Range("O43").Value = 1.7819580999999998
$ws.Range("Q43").Value = 16.6
$ws.Range("S43").Value = 35639.162
$ws.Range("T43").Value = 7.014755285211252
$ws.Range("U43").Value = 15060.240963855422
$ws.Range("J44").Value = 32.7
$ws.Range("L44").Value = 48997.143
$ws.Range("M44").Value = 10.204676627778074
$ws.Range("N44").Value = 15290.51987767584
$ws.Range("Q44").Value = 33.2
$ws.Range("S44").Value = 71278.324
$ws.Range("T44").Value = 7.014755285211252
$ws.Range("U44").Value = 15060.240963855422
$ws.Range("H48").Value = 0.49617645000000005
$ws.Range("J48").Value = 8.35
$ws.Range("L48").Value = 4961.7645
$ws.Range("M48").Value = 25.19265071931568
$ws.Range("N48").Value = 14970.05988023952
$ws.Range("H49").Value = 0.49617645000000005
$ws.Range("J49").Value = 8.35
$ws.Range("L49").Value = 4961.7645
$ws.Range("M49").Value = 25.19265071931568
$ws.Range("N49").Value = 14970.05988023952
$ws.Range("H50").Value = 0.49617645000000005
$ws.Range("J50").Value = 16.7
$ws.Range("L50").Value = 9923.529
$ws.Range("M50").Value = 25.19265071931568
$ws.Range("N50").Value = 14970.05988023952
$ws.Range("J51").Value = 33.4
$ws.Range("L51").Value = 19847.058
$ws.Range("M51").Value = 25.19265071931568
$ws.Range("N51").Value = 14970.05988023952

# ---------------------------------------------------------------------------
# 3. ShrinkToFit formatting for the section-header cells on "Flow Chart"
# ---------------------------------------------------------------------------
$flow = $wb.Worksheets.Item("Flow Chart")
$headerCells = @("B7", "B16", "B25", "B34", "B43", "B52")
foreach ($addr in $headerCells) {
    $r = $flow.Range($addr)
    $r.Font.Name = "Calibri Light"
    $r.Font.Size = 16
    $r.Font.Color = 6378045
    $r.HorizontalAlignment = -4131
    $r.VerticalAlignment = -4108
    $r.ShrinkToFit = $true
}
